# Re-derived transition-probability matrix after adding more simulated games
# (team_specific_matrix / Tulsa_A sheet): counts-per-state were incremented and
# each row was renormalized to counts/rowTotal. Writing the recomputed cell values
# directly reproduces the same renormalization without needing the raw game log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1916167664670659
$ws.Range("C2").Value = 0.5808383233532934
$ws.Range("J2").Value = 0.02095808383233533
$ws.Range("P2").Value = 0.1467065868263473
$ws.Range("S2").Value = 0.05988023952095808
$ws.Range("B3").Value = 0.02010050251256281
$ws.Range("C3").Value = 0.02010050251256281
$ws.Range("J3").Value = 0.04020100502512563
$ws.Range("P3").Value = 0.7437185929648241
$ws.Range("S3").Value = 0.1758793969849246
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.1833333333333333
$ws.Range("B6").Value = 0.06060606060606061
$ws.Range("D6").Value = 0.02525252525252525
$ws.Range("F6").Value = 0.0505050505050505
$ws.Range("J6").Value = 0.2929292929292929
$ws.Range("O6").Value = 0.02525252525252525
$ws.Range("Q6").Value = 0.1565656565656566
$ws.Range("R6").Value = 0.0707070707070707
$ws.Range("S6").Value = 0.3181818181818182
$ws.Range("B7").Value = 0.09042553191489362
$ws.Range("D7").Value = 0.05319148936170213
$ws.Range("F7").Value = 0.05319148936170213
$ws.Range("J7").Value = 0.09042553191489362
$ws.Range("Q7").Value = 0.1595744680851064
$ws.Range("R7").Value = 0.06914893617021277
$ws.Range("S7").Value = 0.4840425531914894
$ws.Range("B8").Value = 0.1100917431192661
$ws.Range("D8").Value = 0.009174311926605505
$ws.Range("F8").Value = 0.04587155963302753
$ws.Range("J8").Value = 0.1032110091743119
$ws.Range("O8").Value = 0.03211009174311927
$ws.Range("Q8").Value = 0.1811926605504587
$ws.Range("R8").Value = 0.07568807339449542
$ws.Range("S8").Value = 0.4426605504587156
$ws.Range("B9").Value = 0.1621621621621622
$ws.Range("D9").Value = 0.02702702702702703
$ws.Range("E9").Value = 0.005405405405405406
$ws.Range("F9").Value = 0.06486486486486487
$ws.Range("J9").Value = 0.08108108108108109
$ws.Range("O9").Value = 0.01081081081081081
$ws.Range("Q9").Value = 0.1351351351351351
$ws.Range("R9").Value = 0.07027027027027027
$ws.Range("S9").Value = 0.4432432432432433
$ws.Range("B10").Value = 0.1230893000804505
$ws.Range("D10").Value = 0.03137570394207562
$ws.Range("F10").Value = 0.05309734513274336
$ws.Range("J10").Value = 0.1037811745776347
$ws.Range("O10").Value = 0.01448109412711183
$ws.Range("Q10").Value = 0.1777956556717619
$ws.Range("R10").Value = 0.09573612228479485
$ws.Range("S10").Value = 0.4006436041834272
$ws.Range("F11").Value = 0.006493506493506494
$ws.Range("G11").Value = 0.1525974025974026
$ws.Range("J11").Value = 0.1136363636363636
$ws.Range("K11").Value = 0.1980519480519481
$ws.Range("L11").Value = 0.512987012987013
$ws.Range("S11").Value = 0.01623376623376623
$ws.Range("G12").Value = 0.7636363636363637
$ws.Range("J12").Value = 0.1393939393939394
$ws.Range("K12").Value = 0.006060606060606061
$ws.Range("L12").Value = 0.04242424242424243
$ws.Range("S12").Value = 0.04848484848484848
$ws.Range("G13").Value = 0.6410256410256411
$ws.Range("J13").Value = 0.282051282051282
$ws.Range("S13").Value = 0.07692307692307693
$ws.Range("F15").Value = 0.01657458563535912
$ws.Range("H15").Value = 0.1215469613259668
$ws.Range("I15").Value = 0.07734806629834254
$ws.Range("J15").Value = 0.3867403314917127
$ws.Range("K15").Value = 0.04972375690607735
$ws.Range("M15").Value = 0.01104972375690608
$ws.Range("O15").Value = 0.03314917127071823
$ws.Range("S15").Value = 0.3038674033149171
$ws.Range("F16").Value = 0.03333333333333333
$ws.Range("H16").Value = 0.1791666666666667
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.4083333333333333
$ws.Range("K16").Value = 0.1208333333333333
$ws.Range("M16").Value = 0.01666666666666667
$ws.Range("O16").Value = 0.02916666666666667
$ws.Range("S16").Value = 0.15
$ws.Range("F17").Value = 0.02349869451697128
$ws.Range("H17").Value = 0.1775456919060052
$ws.Range("I17").Value = 0.0731070496083551
$ws.Range("J17").Value = 0.4099216710182768
$ws.Range("K17").Value = 0.1018276762402089
$ws.Range("M17").Value = 0.02610966057441253
$ws.Range("O17").Value = 0.06527415143603134
$ws.Range("S17").Value = 0.1227154046997389
$ws.Range("F18").Value = 0.015625
$ws.Range("H18").Value = 0.1979166666666667
$ws.Range("I18").Value = 0.07291666666666667
$ws.Range("J18").Value = 0.4322916666666667
$ws.Range("M18").Value = 0.02083333333333333
$ws.Range("O18").Value = 0.0625
$ws.Range("S18").Value = 0.1145833333333333
$ws.Range("F19").Value = 0.02059496567505721
$ws.Range("H19").Value = 0.2051868802440885
$ws.Range("I19").Value = 0.08619374523264683
$ws.Range("J19").Value = 0.3874904652936689
$ws.Range("K19").Value = 0.1083142639206712
$ws.Range("M19").Value = 0.01449275362318841
$ws.Range("N19").Value = 0.0007627765064836003
$ws.Range("O19").Value = 0.05339435545385202
$ws.Range("S19").Value = 0.1235697940503433
